$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows for "第九届环形宇宙动漫游戏嘉年华" and
# "心动恋章·冬日序国乙&代号鸢同人only" events: update "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3087
$ws1.Range("F5").Value = 101

# Sheet "全部类型" (All types) - same two events repeated, update the same column
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3087
$ws4.Range("F10").Value = 101
